# 10/27 gassed up shordyyy sed i need my fordyyyy
#
# The first inline picture (originally named "image1.png", embedding
# media/image1.png) gets resized down, and in the process Word swaps the
# display "name" labels on the two pictures' <wp:docPr>/<pic:cNvPr>
# (the embedded image data / relationship ids are untouched - only the
# name= labels trade places).
#
# InlineShape has no scriptable "Name" property in the Word OM (only
# floating Shape objects expose .Name), so the label swap is done via
# surgical edits of the package's WordOpenXML, while the resize itself
# uses the normal InlineShape.Width/Height API.

$d = $word.ActiveDocument

# --- 1) Resize the first picture (InlineShapes.Item(1)) -------------------
# 5943600 x 4470400 EMU  ->  5624513 x 4227398 EMU (EMU / 12700 = points)
$firstShape = $d.InlineShapes.Item(1)
$firstShape.Width = 442.87503937007875
$firstShape.Height = 332.8659842519685

# --- 2) Swap the two pictures' docPr/cNvPr "name" labels -------------------
function Replace-FirstAfter($str, $search, $replacement, $startAt) {
    $idx = $str.IndexOf($search, $startAt)
    if ($idx -lt 0) {
        throw "pattern not found: $search"
    }
    $before = $str.Substring(0, $idx)
    $after = $str.Substring($idx + $search.Length)
    return $before + $replacement + $after
}

$xml = $d.WordOpenXML

# wp:docPr id="1" (first picture): name image1.png -> image2.png
$pos = $xml.IndexOf('<wp:docPr id="1" name="image1.png"/>')
$xml = Replace-FirstAfter $xml '<wp:docPr id="1" name="image1.png"/>' '<wp:docPr id="1" name="image2.png"/>' 0

# pic:cNvPr that belongs to the same picture comes right after its docPr
$xml = Replace-FirstAfter $xml '<pic:cNvPr id="0" name="image1.png"/>' '<pic:cNvPr id="0" name="image2.png"/>' $pos

# wp:docPr id="2" (second picture): name image2.png -> image1.png
$pos2 = $xml.IndexOf('<wp:docPr id="2" name="image2.png"/>')
$xml = Replace-FirstAfter $xml '<wp:docPr id="2" name="image2.png"/>' '<wp:docPr id="2" name="image1.png"/>' 0

$xml = Replace-FirstAfter $xml '<pic:cNvPr id="0" name="image2.png"/>' '<pic:cNvPr id="0" name="image1.png"/>' $pos2

$d.WordOpenXML = $xml
